$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add the new one right after it -----------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Data Kapal"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Pediksi Akurasi"

# --- Header row (row 1) ------------------------------------------------
$ws2.Range("A1").Value = "No"
$ws2.Range("B1").Value = "Vessel Name"
$ws2.Range("C1").Value = "Vessel Type "
$ws2.Range("D1").Value = "Build in (year)"
$ws2.Range("E1").Value = "Age(year)"
$ws2.Range("F1").Value = "Gross Tonnage"
$ws2.Range("G1").Value = "Summer Deadweight (t)"
$ws2.Range("H1").Value = "Length (m)"
$ws2.Range("I1").Value = "Width (m)"
$ws2.Range("J1").Value = "Accident Risk"

# Alignment first (reuses existing "center" style), fill colour after
# (reuses the existing yellow-fill "center" style for J1 without minting an
# orphan combination).
$ws2.Range("A1:I1").HorizontalAlignment = -4108
$ws2.Range("J1").HorizontalAlignment = -4108
$ws2.Range("J1").Interior.Color = 65535

# --- Data rows -----------------------------------------------------------
$data = @(
  @(1,"TANTO SALAM","Container Ship ",2019,5,10461,11164,137,23,"Low"),
  @(2,"DHT SUNDARBANS","Crude Oil Tanker",2012,12,161513,318123,333,60,"Low"),
  @(3,"SEAWAYS CAPE HENRY","Crude Oil Tanker",2016,8,161319,300932,333,60,"Low"),
  @(4,"NORTHERN RANGER","Passenger/General Cargo Ship",1986,38,2573,662,72,16,"High"),
  @(5,"NORCON GALATEA"," Passenger/General Cargo Ship",1968,56,387,179,41,11,"High"),
  @(6,"PO LIBERTE","Passenger/Ro-Ro Cargo Ship",2023,1,47653,8850,230,31,"Low"),
  @(7,"SEABOURN SOJOURN","Passenger (Cruise) Ship",2010,14,32477,3780,198,26,"High"),
  @(8,"DELTA AMAZON","Crude Oil Tanker",2015,9,166178,319896,333,60,"Low"),
  @(9,"CHI-CHEEMAUN","Passenger/Ro-Ro Cargo",1974,50,6991,855,111,19,"High"),
  @(10,"TORM LAURA","Chemical/Oil Products",2008,16,29300,53160,183,32,"High")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $ws2.Cells.Item($r,1).Value = $row[0]
  $ws2.Cells.Item($r,2).Value = $row[1]
  $ws2.Cells.Item($r,3).Value = $row[2]
  $ws2.Cells.Item($r,4).Value = $row[3]
  $ws2.Cells.Item($r,5).Value = $row[4]
  $ws2.Cells.Item($r,6).Value = $row[5]
  $ws2.Cells.Item($r,7).Value = $row[6]
  $ws2.Cells.Item($r,8).Value = $row[7]
  $ws2.Cells.Item($r,9).Value = $row[8]
  $ws2.Cells.Item($r,10).Value = $row[9]
}

# Row 7 (PO LIBERTE) has its vessel-name cell word-wrapped.
$ws2.Range("B7").WrapText = $true

# --- Column widths ---------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 8.0
$ws2.Columns.Item(2).ColumnWidth = 24.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 25.0
$ws2.Range($ws2.Columns.Item(4), $ws2.Columns.Item(9)).ColumnWidth = 18.333333333333332
$ws2.Columns.Item(10).ColumnWidth = 14.0

# --- sheetView state -------------------------------------------------------
$ws2.Range("G19").Select()

# Sheet1 selection / view reset
$ws1.Range("A1:J1").Select()
$ws2.Activate()
